# Applies the "Update bug fixes for home page icons, filter menu and forms"
# commit to the workbook: adds four new rows (17-20) to the "Bugs" sheet
# documenting new bug fixes, tweaks a few existing row heights, nudges the
# screenshot picture that sits below the bug table, and updates the saved
# selection / active-sheet state to match where the author ended up.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Bugs" sheet - add the four new bug-fix rows (17-20)
# ---------------------------------------------------------------------
$wsBugs = $wb.Worksheets.Item("Bugs")

# Row 16 is the template: same font/fill/border pattern (B/C plain cells,
# D spacer, E:H merged "Outcome/Comments" cell). Clone its formatting down
# into the four new rows, then merge each new E:H block.
$wsBugs.Range("B16:H16").Copy()
$wsBugs.Range("B17:H17").PasteSpecial(-4122)
$wsBugs.Range("B18:H18").PasteSpecial(-4122)
$wsBugs.Range("B19:H19").PasteSpecial(-4122)
$wsBugs.Range("B20:H20").PasteSpecial(-4122)

$wsBugs.Range("E17:H17").Merge()
$wsBugs.Range("E18:H18").Merge()
$wsBugs.Range("E19:H19").Merge()
$wsBugs.Range("E20:H20").Merge()

# Fill in the new content. Order matches the order the strings were
# originally authored in (the "Outcome/Comments" cell for row 17 - the
# EmailJS issue - was written last, after the other three rows).
$wsBugs.Range("B17").Value = "Order Form to send emails via EmailJS"
$wsBugs.Range("C17").Value = '412 Error found in the console: "Preconditioned Fail"'

$wsBugs.Range("B18").Value = "Homepage Icons to be clickable to display company information"
$wsBugs.Range("C18").Value = "First click attempt not functioning and general design not condusive to good UX"
$wsBugs.Range("E18").Value = "Removed JS onClick functionality and replaced with CSS on hover pseudo class"

$wsBugs.Range("B19").Value = "Filter buttons to highlight yellow when active"
$wsBugs.Range("C19").Value = "All products button staying highlighed when others were in use"
$wsBugs.Range("E19").Value = "1. Remove currentBtnClickListner function from being declared as a variable."

$wsBugs.Range("B20").Value = "Forms to handle invalid or empty inputs"
$wsBugs.Range("C20").Value = "Forms accepting empty white space in fields"
$wsBugs.Range("E20").Value = "1.Submission failed despite required attribute present in inputs fields & JS validating against empty or one escape entered. I added an extensive JS function to validate all fields which triggers a modal on submission if an error is present, .trim method used to remove empty spaces."

$wsBugs.Range("E17").Value = "1. Opened emailJS platform and advised by alert to reconnect my email, so deleted it and re-added it. Issue appears to be resolved. Researched & contacted tutor support to find out the cause of the error but they were not familiar with it or how to prevent it in future. Suspect it is caused by the EmailJS platform and outside my code's control. A modal is present to alert the user that their submission failed, I tested it by altering the API Key and it appeared/was called successfully. "

# Row heights - new rows plus re-wrapped existing rows (5, 10, 12)
$wsBugs.Rows.Item(5).RowHeight = 51.75
$wsBugs.Rows.Item(10).RowHeight = 384
$wsBugs.Rows.Item(12).RowHeight = 122.25
$wsBugs.Rows.Item(17).RowHeight = 77.25
$wsBugs.Rows.Item(18).RowHeight = 51.75
$wsBugs.Rows.Item(19).RowHeight = 39
$wsBugs.Rows.Item(20).RowHeight = 53.25

# ---------------------------------------------------------------------
# 2. Nudge the screenshot picture under the bug table (it was dragged
#    slightly left/up and made a little wider).
# ---------------------------------------------------------------------
$shp = $wsBugs.Shapes.Item(1)
$shp.Left = $shp.Left - 1.9644094488188977
$shp.Top = $shp.Top - 24.821417322834645
$shp.Width = $shp.Width + 2.1429133858267715

# ---------------------------------------------------------------------
# 3. Update the saved selection on "Features" (was the active sheet,
#    now isn't) and make "Bugs" the active sheet with its own selection.
# ---------------------------------------------------------------------
$wsFeatures = $wb.Worksheets.Item("Features")
$wsFeatures.Range("D27").Select()

$wsBugs.Activate()
$wsBugs.Range("J10").Select()
